$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Client_Project")

# Fix mismatched Number codes for the two WSNB-WC rows before inserting the new row
$t9 = $ws.Range("C9").Text
$t10 = $ws.Range("C10").Text
$ws.Range("C9").Value = "'" + $t10
$ws.Range("C10").Value = "'" + $t9

# Insert new project time code row
$ws.Rows.Item(9).Insert()
$ws.Range("B9").Value = "ZZZ (Marketing Specific Prospect)"
$ws.Range("C9").Value = "'00000000"
$ws.Range("D9").Value = 101
